$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.071.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.02%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.929.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.19%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "482.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.25%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.64%  "

# Row 7
$ws.Range("E7").Value = "  -2.45%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.722"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.58%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.49%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000356"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +12.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.32%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.18%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.559.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.16%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.957.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.38%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.13%  "

# Row 17
$ws.Range("E17").Value = "  -0.34%  "

# Row 18
$ws.Range("E18").Value = "  -2.51%  "

# Row 19
$ws.Range("E19").Value = "  -3.84%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.107.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.66%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "434.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.59%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.84%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.96%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +14.78%  "

# Row 26
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.73%  "

# Row 27
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.47%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.49%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.50%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "713.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.95%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.129"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.04%  "

# Row 32
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.44%  "

# Row 33
$ws.Range("E33").Value = "  +2.50%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0936"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +35.99%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.43%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.151"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.15%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.88%  "

# Row 39
$ws.Range("E39").Value = "  -0.18%  "

# Row 40
$ws.Range("E40").Value = "  -2.30%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.54%  "

# Row 42
$ws.Range("E42").Value = "  +6.94%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.26%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.338"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.16%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.140"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.50%  "

# Row 46
$ws.Range("E46").Value = "  -0.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.08%  "

# Row 48
$ws.Range("E48").Value = "  -1.44%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "148.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.77%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.45%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.64%  "
